$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("01_IB전략컨설팅부")

# Delete the row for "미래에셋비전스팩4호" (수요예측 시작일 2024-05-13), which is row 20.
# This shifts the following row ("신한글로벌액티브리츠") up to become row 20.
$ws.Rows.Item(20).Delete()
